# Apply the changes described by the diff to both the "展览" (sheet 1)
# and "全部类型" (sheet 4) worksheets. Sheet 4 contains the same
# exhibition rows as sheet 1 (plus some rows from other sheets mixed
# in, sorted by date), so the row numbers are not a uniform offset of
# sheet 1's row numbers - an explicit per-row mapping is used below.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must stay plain text even when it looks
# like a date (e.g. "2024-09-16"), which Excel would otherwise silently
# convert into a date serial number. Forcing the NumberFormat to "@"
# (Text) before assigning the value keeps it literal; ClearFormats
# afterwards removes the now-unneeded text format from the cell again
# so its formatting matches the original (unformatted) cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Update-ExhibitionSheet($ws, $r3, $r5, $r6, $r10, $r13, $r14, $r15, $r17, $r19, $r25, $r26) {

    # --- simple "想去人数" (F column) bumps ---
    $ws.Cells.Item($r3, 6).Value = 12934   # was 12926
    $ws.Cells.Item($r10, 6).Value = 12893  # was 12880
    $ws.Cells.Item($r13, 6).Value = 8696   # was 8694
    $ws.Cells.Item($r14, 6).Value = 7700   # was 7696
    $ws.Cells.Item($r15, 6).Value = 200    # was 198
    $ws.Cells.Item($r19, 6).Value = 985    # was 984
    $ws.Cells.Item($r25, 6).Value = 88     # was 87
    $ws.Cells.Item($r26, 6).Value = 5220   # was 5219

    # --- G column on the "明日方舟ONLY#2024~佑桑柔" row becomes "不可售" ---
    $ws.Cells.Item($r17, 7).Value = "不可售"

    # --- rows 5 & 6 (AME嘉年华 / 明日方舟同人展ONLY) swap their content ---

    # New row 5 content (was the AME row, becomes 明日方舟同人展ONLY)
    Set-TextValue $ws.Cells.Item($r5, 2) "2024-09-16"
    $ws.Cells.Item($r5, 3).Value = "苏州·明日方舟同人展ONLY"
    $ws.Cells.Item($r5, 4).Value = "小外滩街苏州小外滩 元和塘美术馆"
    $ws.Cells.Item($r5, 5).Value = "2024.09.16 09:00-09.17 16:00"
    $ws.Cells.Item($r5, 6).Value = 77
    $ws.Cells.Item($r5, 7).Value = 69
    $ws.Cells.Item($r5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90913"
    $ws.Cells.Item($r5, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/0wvOxXfo1724120638374.png"

    # New row 6 content (was the 明日方舟同人展ONLY row, becomes AME嘉年华,
    # note: B column (date) stays "2024-09-16" and E column gets an
    # updated time range rather than the old row-5 value)
    Set-TextValue $ws.Cells.Item($r6, 2) "2024-09-16"
    $ws.Cells.Item($r6, 3).Value = "苏州·第二届-AME动漫嘉年华（免费展）"
    $ws.Cells.Item($r6, 4).Value = "东吴南路179号 龙湖苏州东吴天街"
    $ws.Cells.Item($r6, 5).Value = "2024.09.16 10:00-09.16 13:00"
    $ws.Cells.Item($r6, 6).Value = 89
    $ws.Cells.Item($r6, 7).Value = 39.9
    $ws.Cells.Item($r6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90288"
    $ws.Cells.Item($r6, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/e3uZDian1722619198829.png"
}

# Sheet 1 = "展览" -- rows map 1:1 to themselves
$wsExpo = $wb.Worksheets.Item(1)
Update-ExhibitionSheet $wsExpo 3 5 6 10 13 14 15 17 19 25 26

# Sheet 4 = "全部类型" -- same logical rows, but at different row numbers
# (rows 1-21 of sheet1 sit one row lower here; rows 22-26 sit two rows
# lower, because an extra "演出" row is interleaved at sheet4 row 23).
$wsAll = $wb.Worksheets.Item(4)
Update-ExhibitionSheet $wsAll 4 6 7 11 14 15 16 18 20 27 28
